$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in the data rows
$ws.Range("A2").Value = "TestAA"
$ws.Range("B2").Value = 9873787678
$ws.Range("D2").Value = "Juli Co. Ltd"
$ws.Range("B3").Value = 9858787678
$ws.Range("B4").Value = 9878287678
$ws.Range("B5").Value = 9878792678

# Update the selected cell
$ws.Range("D8").Select()
